# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Brayan Yesid Hernandez Argel" record (previously the last row, r44)
# moves up to become the first worker row (r16). The existing rows shift
# down by one, and Jaime Luis Rodriguez Parra's 28 monthly "Periodo Mora"
# rows are re-ordered from descending (2010 -> 1808) to ascending
# (1808 -> 2010), keeping each period's own Salario Basico / Valor Mora
# pair intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 16-44: Tipo Doc, N Doc, Nombre, Periodo, Salario Basico (F), Valor Mora (G)
$rows = @(
    @(16, "CC", "1143396690", "BRAYAN YESID HERNANDEZ ARGEL", "1610", 27578, 689454),
    @(17, "CC", "1143367576", "CRISTIAN PEREZ PESTANA", "1701", 27578, 689455),
    @(18, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1808", 31249, 781242),
    @(19, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1809", 31249, 781242),
    @(20, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1810", 31249, 781242),
    @(21, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1811", 31249, 781242),
    @(22, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1812", 31249, 781242),
    @(23, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1901", 31249, 781242),
    @(24, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1902", 31249, 781242),
    @(25, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1903", 31249, 781242),
    @(26, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1904", 31249, 781242),
    @(27, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1905", 31249, 781242),
    @(28, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1906", 31249, 781242),
    @(29, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1907", 31249, 781242),
    @(30, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1908", 31249, 781242),
    @(31, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1909", 31249, 781242),
    @(32, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1910", 31249, 781242),
    @(33, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1911", 31249, 781242),
    @(34, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "1912", 31249, 781242),
    @(35, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2001", 31249, 781242),
    @(36, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2002", 31249, 781242),
    @(37, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2003", 31249, 781242),
    @(38, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2004", 31249, 781242),
    @(39, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2005", 31249, 781242),
    @(40, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2006", 31249, 781242),
    @(41, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2007", 31249, 781242),
    @(42, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2008", 31249, 781242),
    @(43, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2009", 31249, 781242),
    @(44, "CC", "9295299", "JAIME LUIS RODRIGUEZ PARRA", "2010", 26041, 781242)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
